$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A116").Value = "20180404_02_002_004"
$ws.Range("A118").Value = "20180405_01_004_006"
$ws.Range("A124").Value = "20180405_01_078_082"
$ws.Range("A127").Value = "20180405_01_134_136"
$ws.Range("A129").Value = "20180418_01_001_003"
$ws.Range("A132").Value = "20180419_02_035_036"
$ws.Range("A133").Value = "20180419_02_040_043"
$ws.Range("A137").Value = "20180419_02_123_126"
$ws.Range("A144").Value = "20180423_01_093_094"
$ws.Range("A145").Value = "20180423_01_097_100"
$ws.Range("A151").Value = "20180426_01_016_018"
$ws.Range("A158").Value = "20180430_01_030_032"

$ws.Range("A158").Select()
try { $excel.ActiveWindow.ScrollRow = 151 } catch { }
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }
